$d = $word.ActiveDocument

# Helper: replace an exact run of text that lives inside a specific
# paragraph (identified by its 1-based paragraph index) with new text,
# optionally marking the new run bold/italic. Using Range.InsertXML (scoped
# to a freshly-created Range - NOT the Find range itself, and NOT the whole
# story range) swaps only the matched text's run while leaving any sibling
# runs (such as the stray empty <w:r/> that precedes many paragraphs in this
# document) completely untouched.
function Replace-TextInParagraph($ParaIndex, $OldText, $NewText, $Bold, $Italic) {
    $p = $d.Paragraphs($ParaIndex)
    $searchRng = $d.Range($p.Range.Start, $p.Range.End)
    $found = $searchRng.Find.Execute($OldText, $true, $false, $false, $false, $false, $true)
    if (-not $found) {
        Write-Host "NOT FOUND in paragraph" $ParaIndex ":" $OldText
        return
    }

    $targetRng = $d.Range($searchRng.Start, $searchRng.End)

    $rPr = ""
    if ($Bold)   { $rPr += "<w:b/>" }
    if ($Italic) { $rPr += "<w:i/>" }

    $runXml = "<w:r>"
    if ($rPr -ne "") { $runXml += "<w:rPr>$rPr</w:rPr>" }
    $runXml += "<w:t>$NewText</w:t></w:r>"

    $fragment = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        "<w:body><w:p>$runXml</w:p></w:body>" +
        '</w:document></pkg:xmlData></pkg:part></pkg:package>'

    $targetRng.InsertXML($fragment)
}

# 1. Title (Heading1, paragraph 1) - no explicit run formatting.
Replace-TextInParagraph 1 `
    "Play Ganesha Jr Free: Impressive Graphics & Exciting Gameplay" `
    "Play Ganesha Jr for Free - Exciting Slot Game with Impressive Graphics" `
    $false $false

# 2. "What we like" bullet list items (paragraphs 42-45).
Replace-TextInParagraph 42 `
    "Impressive graphics and sound effects" `
    "High-quality graphics and representations" `
    $false $false

Replace-TextInParagraph 43 `
    "Numerous high-value symbols for great rewards" `
    "Satisfying sound effects" `
    $false $false

Replace-TextInParagraph 44 `
    "50 free spins from the start" `
    "Numerous opportunities to win with special symbols" `
    $false $false

Replace-TextInParagraph 45 `
    "Entertaining and well-designed" `
    "50 free spins for players" `
    $false $false

# 3. "What we don't like" bullet list items (paragraphs 47-48).
Replace-TextInParagraph 47 `
    "No progressive jackpot" `
    "Limited information on specific gameplay features" `
    $false $false

Replace-TextInParagraph 48 `
    "No bonus game" `
    "No progressive jackpot" `
    $false $false

# 4. Bold restatement of the title (paragraph 49) - explicit <w:b/>.
Replace-TextInParagraph 49 `
    "Play Ganesha Jr Free: Impressive Graphics & Exciting Gameplay" `
    "Play Ganesha Jr for Free - Exciting Slot Game with Impressive Graphics" `
    $true $false

# 5. Italic meta description (paragraph 50) - explicit <w:i/>.
Replace-TextInParagraph 50 `
    "Read our review of Ganesha Jr, an exciting and well-designed online slot game by CQ9 Gaming with impressive graphics and numerous opportunities to win. Play for free!" `
    "Play Ganesha Jr for free and enjoy high-quality graphics and numerous opportunities to win rewards." `
    $false $true
